# Añadido el reloj: se agregan los segundos ("ss") al display del reloj en
# la tabla de la Slide 1 (Menú Inicial), separados de los minutos por ":".
# Antes las celdas 6, 7 y 8 de la primera fila de la tabla estaban vacías
# (solo tenian el endParaRPr); ahora llevan el texto ":", "s" y "s".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tbl = $sh.Table

$tbl.Cell(1, 6).Shape.TextFrame.TextRange.Text = ":"
$tbl.Cell(1, 7).Shape.TextFrame.TextRange.Text = "s"
$tbl.Cell(1, 8).Shape.TextFrame.TextRange.Text = "s"
